
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "C" column (the old C, with the
# per-analyst rating detail / "UN" placeholder, shifts right to E).
$ws.Columns("C:D").Insert()

# New header row: newest date columns are prepended, so B/C/D/E now read
# Jun_17, Jun_15, Jun_13, Jun_10 left-to-right.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# The two freshly inserted columns get the same "UN" placeholder used
# throughout column B for every data row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Match column C's original width (8 characters) on the two new columns and
# keep it on the shifted-over column as well.
$ws.Range("C1:E1").ColumnWidth = 7.166666666666667
